$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.0436123049819568
$ws.Range("C2").Value = 0.019083653967749
$ws.Range("D2").Value = 0.0681409559961647
$ws.Range("B3").Value = 0.0531128486474015
$ws.Range("C3").Value = 0.017971696962778
$ws.Range("D3").Value = 0.088254000332025
$ws.Range("B4").Value = 0.00119832579274969
$ws.Range("C4").Value = -0.0787656735741821
$ws.Range("D4").Value = 0.0811623251596815
$ws.Range("B5").Value = 0.0540511173859778
$ws.Range("C5").Value = -0.0209004811205222
$ws.Range("D5").Value = 0.129002715892478
$ws.Range("B6").Value = 0.0109083790390346
$ws.Range("C6").Value = -0.078808524563807
$ws.Range("D6").Value = 0.100625282641876
$ws.Range("B7").Value = 0.0529192867753305
$ws.Range("C7").Value = -0.0547675567118714
$ws.Range("D7").Value = 0.160606130262532
$ws.Range("B8").Value = 0.129460428672861
$ws.Range("C8").Value = 0.0297698628802134
$ws.Range("D8").Value = 0.229150994465509
$ws.Range("B9").Value = 0.0808082726394657
$ws.Range("C9").Value = -0.00592646017538637
$ws.Range("D9").Value = 0.167543005454318
$ws.Range("B10").Value = 0.0806919758891738
$ws.Range("C10").Value = -0.0405955523514283
$ws.Range("D10").Value = 0.201979504129776
$ws.Range("B11").Value = 0.040516322743288
$ws.Range("C11").Value = -0.0119826936379433
$ws.Range("D11").Value = 0.0930153391245193
$ws.Range("B12").Value = 0.036707291653055
$ws.Range("C12").Value = -0.00601192427354733
$ws.Range("D12").Value = 0.0794265075796573
